# Add a new "Player Info" worksheet before the existing "ODI Batting" sheet,
# and update the "ODI Batting" sheet's MATCH_CARD_LINK column to MATCH_CODE
# (storing just the numeric match code instead of the full URL).

$wb = $excel.ActiveWorkbook
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 1. Create the new "Player Info" sheet, placed before "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").Borders.LineStyle = 1
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160

$playerInfo.Range("A2").Value = "'5829"
$playerInfo.Range("B2").Value = "Kyle Verreynne"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- 2. Update "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE,
#        and replace full URLs with just the trailing match code number ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$usedRange = $battingSheet.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split "MatchCode="
        if ($parts.Length -gt 1) {
            $cell.Value = "'" + $parts[1]
        }
    }
}
